# Auto-generated edit script applying the Seraph_Profits.xlsx diff
# Updates numeric value cells (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1349.3334
$ws.Range("I17").Value = 999
$ws.Range("K17").Value = 2997
$ws.Range("M17").Value = -2829
# Row 40
$ws.Range("H40").Value = 2118.476
$ws.Range("J40").Value = 2372.6365
$ws.Range("L40").Value = 2372.6365
$ws.Range("N40").Value = -2722.6365
# Row 88
$ws.Range("H88").Value = 2421.6
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 2827
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 2827
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -3639
# Row 91
$ws.Range("H91").Value = 2421.6
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 2827
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 2827
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -5635
# Row 129
$ws.Range("H129").Value = 2250.818
$ws.Range("I129").Value = 1747.125
$ws.Range("J129").Value = 3594
$ws.Range("K129").Value = 5241.375
$ws.Range("L129").Value = 10782
$ws.Range("M129").Value = -241.375
$ws.Range("N129").Value = -20782
# Row 135
$ws.Range("H135").Value = 1020.31036
$ws.Range("I135").Value = 590.4761999999999
$ws.Range("J135").Value = 2148.625
$ws.Range("K135").Value = 5314.2858
$ws.Range("L135").Value = 19337.625
$ws.Range("M135").Value = -2779.2858
$ws.Range("N135").Value = -24407.625

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 2617.6
$ws.Range("I102").Value = 2617.6
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2617.6
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -995.5999999999999
$ws.Range("N102").ClearContents()
# Row 132
$ws.Range("H132").Value = 2434.5715
$ws.Range("I132").Value = 1849.04
$ws.Range("K132").Value = 5547.12
$ws.Range("M132").Value = -3017.12

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 56666.668
$ws.Range("J40").Value = 56666.668
$ws.Range("L40").Value = 56666.668
$ws.Range("N40").Value = -57196.668
# Row 43
$ws.Range("H43").Value = 80000
$ws.Range("J43").Value = 80000
$ws.Range("L43").Value = 80000
$ws.Range("N43").Value = -80362
# Row 87
$ws.Range("H87").Value = 91666.664
$ws.Range("I87").Value = 75000
$ws.Range("J87").Value = 100000
$ws.Range("K87").Value = 75000
$ws.Range("L87").Value = 100000
$ws.Range("M87").Value = -73752
$ws.Range("N87").Value = -102496
# Row 90
$ws.Range("H90").Value = 91666.664
$ws.Range("I90").Value = 75000
$ws.Range("J90").Value = 100000
$ws.Range("K90").Value = 225000
$ws.Range("L90").Value = 300000
$ws.Range("M90").Value = -218760
$ws.Range("N90").Value = -312480
# Row 94
$ws.Range("H94").Value = 720.5263
$ws.Range("I94").Value = 720.5263
$ws.Range("K94").Value = 720.5263
$ws.Range("M94").Value = -269.5263
# Row 96
$ws.Range("H96").Value = 41999.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 41999.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 41999.5
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -47491.5
# Row 99
$ws.Range("H99").Value = 1542.1786
$ws.Range("I99").Value = 1367.24
$ws.Range("K99").Value = 1367.24
$ws.Range("M99").Value = 130.76
# Row 105
$ws.Range("H105").Value = 3494.75
$ws.Range("I105").Value = 2990
$ws.Range("K105").Value = 2990
$ws.Range("M105").Value = -1243
# Row 134
$ws.Range("H134").Value = 2463.1538
$ws.Range("I134").Value = 1702.1
$ws.Range("K134").Value = 5106.299999999999
$ws.Range("M134").Value = -2571.299999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5319.1924
$ws.Range("J31").Value = 7575.1665
$ws.Range("L31").Value = 7575.1665
$ws.Range("N31").Value = -8165.1665
# Row 34
$ws.Range("H34").Value = 5319.1924
$ws.Range("J34").Value = 7575.1665
$ws.Range("L34").Value = 7575.1665
$ws.Range("N34").Value = -7979.1665
# Row 41
$ws.Range("H41").Value = 22692.857
$ws.Range("I41").Value = 5875
$ws.Range("K41").Value = 5875
$ws.Range("M41").Value = -5447
# Row 134
$ws.Range("H134").Value = 2342.658
$ws.Range("I134").Value = 2180.5356
$ws.Range("J134").Value = 2796.6
$ws.Range("K134").Value = 6541.6068
$ws.Range("L134").Value = 8389.799999999999
$ws.Range("M134").Value = -4006.6068
$ws.Range("N134").Value = -13459.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 300
$ws.Range("M13").Value = -132
$ws.Range("N13").Value = -636
# Row 17
$ws.Range("H17").Value = 727.64703
$ws.Range("I17").Value = 93.25
$ws.Range("J17").Value = 2250.2
$ws.Range("K17").Value = 279.75
$ws.Range("L17").Value = 6750.599999999999
$ws.Range("M17").Value = -110.75
$ws.Range("N17").Value = -7088.599999999999
# Row 33
$ws.Range("H33").Value = 90961.63
$ws.Range("J33").Value = 142929
$ws.Range("L33").Value = 857574
$ws.Range("N33").Value = -858140
# Row 46
$ws.Range("H46").Value = 5000299.5
$ws.Range("I46").Value = 599
$ws.Range("K46").Value = 1797
$ws.Range("M46").Value = -1706
# Row 111
$ws.Range("H111").Value = 3245.6
$ws.Range("J111").Value = 3245.6
$ws.Range("L111").Value = 9736.799999999999
$ws.Range("N111").Value = -15870.8
# Row 137
$ws.Range("H137").Value = 4609
$ws.Range("J137").Value = 4666.4
$ws.Range("L137").Value = 13999.2
$ws.Range("N137").Value = -24199.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 38888.5
$ws.Range("J20").Value = 38888.5
$ws.Range("L20").Value = 38888.5
$ws.Range("N20").Value = -39378.5
# Row 24
$ws.Range("H24").Value = 37777.4
$ws.Range("J24").Value = 37777.4
$ws.Range("L24").Value = 37777.4
$ws.Range("N24").Value = -38123.4
# Row 70
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 3333.3333
$ws.Range("J70").Value = 6222.222
$ws.Range("K70").Value = 3333.3333
$ws.Range("L70").Value = 6222.222
$ws.Range("M70").Value = -3063.3333
$ws.Range("N70").Value = -6762.222
# Row 73
$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 3333.3333
$ws.Range("J73").Value = 6222.222
$ws.Range("K73").Value = 3333.3333
$ws.Range("L73").Value = 6222.222
$ws.Range("M73").Value = -2397.3333
$ws.Range("N73").Value = -8094.222
# Row 132
$ws.Range("H132").Value = 2820.8235
$ws.Range("I132").Value = 2182.182
$ws.Range("K132").Value = 6546.545999999999
$ws.Range("M132").Value = -4016.545999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 28653.2
$ws.Range("J20").Value = 27503
$ws.Range("L20").Value = 27503
$ws.Range("N20").Value = -27955
# Row 22
$ws.Range("H22").Value = 6526.4287
$ws.Range("I22").Value = 4000.5
$ws.Range("K22").Value = 4000.5
$ws.Range("M22").Value = -3705.5
# Row 27
$ws.Range("H27").Value = 6526.4287
$ws.Range("I27").Value = 4000.5
$ws.Range("K27").Value = 4000.5
$ws.Range("M27").Value = -3893.5
# Row 55
$ws.Range("H55").Value = 522.1177
$ws.Range("J55").Value = 997.4286
$ws.Range("L55").Value = 997.4286
$ws.Range("N55").Value = -1343.4286
# Row 61
$ws.Range("H61").Value = 3031.3333
$ws.Range("I61").Value = 2469
$ws.Range("K61").Value = 2469
$ws.Range("M61").Value = -2267
# Row 113
$ws.Range("H113").Value = 3031.3333
$ws.Range("I113").Value = 2469
$ws.Range("K113").Value = 2469
$ws.Range("M113").Value = -299
# Row 132
$ws.Range("H132").Value = 2895
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 5273.857
$ws.Range("I136").Value = 5334.0557
$ws.Range("J136").Value = 4912.6665
$ws.Range("K136").Value = 16002.1671
$ws.Range("L136").Value = 14737.9995
$ws.Range("M136").Value = -13452.1671
$ws.Range("N136").Value = -19837.9995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1627.5883
$ws.Range("I126").Value = 1634.2858
$ws.Range("J126").Value = 1634.2858
$ws.Range("K126").Value = 4902.857400000001
$ws.Range("M126").Value = -2432.857400000001
